$d = $word.ActiveDocument
$r1 = $d.Range(0,7)
$r1.Text = "Lesson "
$r2 = $d.Range(7,9)
$r2.Text = "9"
Write-Host $d.Content.Text
